$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$ws = $wb.Worksheets.Item("P_valores")

$ws.Range("C2").Value = 0.5918383337234023
$ws.Range("D2").Value = 0.8732392521185299
$ws.Range("E2").Value = 0.8585639443434183
$ws.Range("F2").Value = 0.2375827175323819

$ws.Range("B3").Value = 0.5918383337234023
$ws.Range("D3").Value = 0.6549773821717224
$ws.Range("E3").Value = 0.2804859867908815
$ws.Range("F3").Value = 0.5583637151360441

$ws.Range("B4").Value = 0.8732392521185299
$ws.Range("C4").Value = 0.6549773821717224
$ws.Range("E4").Value = 0.6828202369083396
$ws.Range("F4").Value = 0.3832568638751579

$ws.Range("B5").Value = 0.8585639443434183
$ws.Range("C5").Value = 0.2804859867908815
$ws.Range("D5").Value = 0.6828202369083396
$ws.Range("F5").Value = 0.1719394901564206

$ws.Range("B6").Value = 0.2375827175323819
$ws.Range("C6").Value = 0.5583637151360441
$ws.Range("D6").Value = 0.3832568638751579
$ws.Range("E6").Value = 0.1719394901564206

# --- Sheet: Estadisticos_DM ---
$ws2 = $wb.Worksheets.Item("Estadisticos_DM")

$ws2.Range("C2").Value = -0.5487248654433636
$ws2.Range("D2").Value = -0.1624935533370932
$ws2.Range("E2").Value = 0.1815165397954219
$ws2.Range("F2").Value = -1.233819633964347

$ws2.Range("B3").Value = 0.5487248654433636
$ws2.Range("D3").Value = 0.4565784363014768
$ws2.Range("E3").Value = 1.122613201377622
$ws2.Range("F3").Value = -0.5995831086647904

$ws2.Range("B4").Value = 0.1624935533370932
$ws2.Range("C4").Value = -0.4565784363014768
$ws2.Range("E4").Value = 0.4172543568248846
$ws2.Range("F4").Value = -0.9001624317704398

$ws2.Range("B5").Value = -0.1815165397954219
$ws2.Range("C5").Value = -1.122613201377622
$ws2.Range("D5").Value = -0.4172543568248846
$ws2.Range("F5").Value = -1.439691460963566

$ws2.Range("B6").Value = 1.233819633964347
$ws2.Range("C6").Value = 0.5995831086647904
$ws2.Range("D6").Value = 0.9001624317704398
$ws2.Range("E6").Value = 1.439691460963566
